$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = -7.948399999999995
$ws.Range("B7").Value = 5.128999999999999
$ws.Range("E7").Value = 16.2023
$ws.Range("A9").Value = -21.9242
$ws.Range("E10").Value = 16.5386
$ws.Range("B12").Value = 5.345099999999997
$ws.Range("A13").Value = -22.24679999999999
$ws.Range("E13").Value = 16.6734
$ws.Range("B14").Value = 6.100500000000003
$ws.Range("D15").Value = -8.841499999999998
$ws.Range("A16").Value = -21.923
$ws.Range("E16").Value = 15.8353
$ws.Range("A18").Value = -22.17460000000001
$ws.Range("B19").Value = 8.8284
$ws.Range("A20").Value = -19.93569999999999
$ws.Range("E20").Value = 15.8903
$ws.Range("E24").Value = 16.43829999999999
$ws.Range("A26").Value = -21.10959999999998
$ws.Range("B26").Value = 4.534400000000004
$ws.Range("A27").Value = -21.55779999999996
$ws.Range("B27").Value = 5.372000000000005
$ws.Range("D28").Value = -8.394199999999998
$ws.Range("A29").Value = -21.73789999999999
$ws.Range("B29").Value = 5.4955
$ws.Range("E32").Value = 16.28379999999998
$ws.Range("D33").Value = -7.474100000000002
$ws.Range("A35").Value = -19.6535
$ws.Range("D35").Value = -9.03709999999999
$ws.Range("A36").Value = -19.98859999999998
$ws.Range("B37").Value = 9.149300000000006
$ws.Range("B38").Value = 4.549699999999999
$ws.Range("D38").Value = -9.123599999999993
$ws.Range("E39").Value = 16.32739999999999
$ws.Range("D43").Value = -8.457099999999999
$ws.Range("D44").Value = -7.250399999999998
$ws.Range("A45").Value = -21.56899999999999
$ws.Range("D45").Value = -7.863099999999999
$ws.Range("B47").Value = 6.280000000000005
$ws.Range("D47").Value = -7.504000000000001
$ws.Range("E47").Value = 17.3775
$ws.Range("E48").Value = 17.5158
$ws.Range("B51").Value = 5.655799999999999
$ws.Range("D51").Value = -7.187499999999996
$ws.Range("B52").Value = 5.288399999999998
$ws.Range("E52").Value = 17.16210000000001
$ws.Range("D54").Value = -8.275499999999997
$ws.Range("A55").Value = -22.21890000000001
$ws.Range("B55").Value = 5.196099999999998
$ws.Range("E56").Value = 16.86460000000001
$ws.Range("A57").Value = -22.0913
$ws.Range("D57").Value = -8.316400000000002
$ws.Range("D62").Value = -8.550799999999999
$ws.Range("D63").Value = -8.347499999999995
$ws.Range("D67").Value = -5.961199999999996
$ws.Range("A69").Value = -21.6664
$ws.Range("B69").Value = 5.557599999999996
$ws.Range("B70").Value = 5.692800000000001
$ws.Range("D70").Value = -6.969100000000002
$ws.Range("A76").Value = -22.02969999999999
$ws.Range("B76").Value = 4.855799999999996
$ws.Range("A78").Value = -19.74269999999998
$ws.Range("B81").Value = 5.697700000000006
$ws.Range("D81").Value = -7.053599999999997
$ws.Range("A82").Value = -22.0798
$ws.Range("A83").Value = -21.92670000000001
$ws.Range("B83").Value = 6.404400000000002
$ws.Range("E84").Value = 16.4239
$ws.Range("D88").Value = -7.274599999999993
$ws.Range("A93").Value = -20.36119999999998
$ws.Range("B94").Value = 5.3381
$ws.Range("D96").Value = -8.342600000000003
$ws.Range("A97").Value = -21.8818
$ws.Range("D99").Value = -8.046899999999994
$ws.Range("B100").Value = 5.8043
$ws.Range("E100").Value = 16.36040000000001
$ws.Range("E101").Value = 16.88790000000002
$ws.Range("B102").Value = 9.05080000000001
